$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "review_topics"
$ws.Range("AB1").Value = "trait_frame"

$ws.Range("N2").Value = "skillassessment; selfreport"

$ws.Range("N3").Value = "skilltraitdifferences; incrementalvalidity; academicachievement;"
$ws.Range("O3").Value = "both"

$ws.Range("N4").Value = "skillassessment; selfreport; skilltraitdifferences; incrementalvalidity; otherreport"

$ws.Range("N5").Value = "ses; gender; age"

$ws.Range("N6").Value = "academicachievement; learningdisabilities"
$ws.Range("O6").Value = "both"

$ws.Range("N7").Value = "academicachievement; learningdisabilities; disabilities"
$ws.Range("O7").Value = "both"

$ws.Range("N8").Value = "nomoligicalnet"
$ws.Range("O8").Value = "both"

$ws.Range("N9").Value = "nomoligicalnet"
$ws.Range("O9").Value = "both"

$ws.Range("N10").Value = "skillassessment; selfreport; shortversions"

$ws.Range("N11").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O11").Value = "both"

$ws.Range("N12").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O12").Value = "both"

$ws.Range("N13").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O13").Value = "both"

$ws.Range("N14").Value = "academicachievement; incrementalvalidity; skilltraitdifferences"
$ws.Range("O14").Value = "both"

$ws.Range("N15").Value = "normativechange; age; gender"

$ws.Range("N16").Value = "academicachievement; incrementalvalidity; skilltraitdifferences"
$ws.Range("O16").Value = "both"

$ws.Range("N17").Value = "volunteering"
$ws.Range("O17").Value = "both"

$ws.Range("N18").Value = "academicachievement; incrementalvalidity; skilltraitdifferences"
$ws.Range("O18").Value = "both"

$ws.Range("N19").Value = "skillassessment; behavioralassessment"
$ws.Range("O19").Value = "both"

$ws.Range("N20").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O20").Value = "both"

$ws.Range("N21").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O21").Value = "both"

$ws.Range("N22").Value = "skillassessment; selfreport; nomologicalnet; languageadaptation"
$ws.Range("O22").Value = "both"

$ws.Range("N23").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O23").Value = "both"

$ws.Range("N24").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O24").Value = "both"

$ws.Range("N25").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O25").Value = "both"

$ws.Range("N26").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O26").Value = "both"

$ws.Range("N27").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O27").Value = "both"

$ws.Range("N28").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O28").Value = "both"

$ws.Range("N29").Value = "skillassessment; selfreport; otherreport; incrementalvalidity; nomologicalnet; skilltraitdifferences"
$ws.Range("O29").Value = "both"

$ws.Range("N30").Value = "theorydevelopment"
